# evo: new merge fields
#
# The "Special styles" merge-field cheat-sheet sheet gets its VAR column (B)
# updated: several placeholders were renamed and a handful of brand new
# merge fields were introduced. Column A (the human-readable NAME) is left
# untouched, exactly like the authored change.
#
# New / renamed fields are written first, in the order they were newly
# introduced by the author, so the workbook's shared-string table grows in
# the same order as the target edit; the remaining rows are then written to
# restore the fields that simply moved around.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- brand new / renamed merge fields (first use order) -------------------
$ws.Range("B27").Value = "[datetime.date]"
$ws.Range("B17").Value = "[res_letterbox.nature_id]"
$ws.Range("B15").Value = "[res_letterbox.type_label]"
$ws.Range("B24").Value = "[res_letterbox.alt_identifier]"
$ws.Range("B6").Value  = "[contact.address_num]"
$ws.Range("B7").Value  = "[contact.address_street]"
$ws.Range("B8").Value  = "[contact.address_complement]"
$ws.Range("B9").Value  = "[contact.address_town]"
$ws.Range("B10").Value = "[contact.address_postal_code]"
$ws.Range("B11").Value = "[contact.address_contry]"

# --- fields that only shifted position, value unchanged -------------------
$ws.Range("B12").Value = "[user.lastname]"
$ws.Range("B13").Value = "[user.firstname]"
$ws.Range("B14").Value = "[res_letterbox.destination]"
$ws.Range("B16").Value = "[res_letterbox.category_id]"
$ws.Range("B18").Value = "[res_letterbox.admission_date]"
$ws.Range("B19").Value = "[res_letterbox.doc_date]"
$ws.Range("B20").Value = "[res_letterbox.process_limit_date]"
$ws.Range("B21").Value = "[res_letterbox.process_notes]"
$ws.Range("B22").Value = "[res_letterbox.closing_date]"
$ws.Range("B23").Value = "[res_letterbox.subject]"
$ws.Range("B25").Value = "[res_letterbox.author]"
$ws.Range("B26").Value = "[res_letterbox.creation_date]"
$ws.Range("B28").Value = "[user.lastname]"
$ws.Range("B29").Value = "[user.firstname]"
$ws.Range("B30").Value = "[user.phone]"
$ws.Range("B31").Value = "[user.mail]"

# Match the author's new active-cell selection on the sheet.
$ws.Range("E19").Select()
